# debugging T+R /= 1
# Update the RGF input sheet: rows 15 & 16 go from cell-type "x" to "o",
# and the three region rows (14-16) get updated Width/Length/Gap values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 and 16: cell type changes from "x" to "o"
$ws.Range("A15").Value = "o"
$ws.Range("A16").Value = "o"

# Width(# of sub unit cell): E14:E16 1 -> 2
$ws.Range("E14").Value = 2
$ws.Range("E15").Value = 2
$ws.Range("E16").Value = 2

# Length(# of unit cell): F14 10 -> 300, F15/F16 0 -> 300
$ws.Range("F14").Value = 300
$ws.Range("F15").Value = 300
$ws.Range("F16").Value = 300

# Gap Open (eV): I14:I16 0 -> 0.05
$ws.Range("I14").Value = 0.05
$ws.Range("I15").Value = 0.05
$ws.Range("I16").Value = 0.05

# Recalculate dependent formulas in B7 (Max ribbon width) and B8 (Max ribbon length)
$excel.Calculate()

# Update the active selection to match the saved view state
$ws.Range("E17").Select()
